$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A updates (rows 5-21)
$ws.Range("A5").Value = "领益智造"
$ws.Range("A6").Value = "利欧股份"
$ws.Range("A7").Value = "北方稀土"
$ws.Range("A8").Value = "岩山科技"
$ws.Range("A9").Value = "天融信"
$ws.Range("A10").Value = "吉视传媒"
$ws.Range("A11").Value = "剑桥科技"
$ws.Range("A12").Value = "合力泰"
$ws.Range("A13").Value = "新易盛"
$ws.Range("A14").Value = "启明信息"
$ws.Range("A15").Value = "华银电力"
$ws.Range("A16").Value = "英维克"
$ws.Range("A17").Value = "瑞芯微"
$ws.Range("A18").Value = "步步高"
$ws.Range("A19").Value = "方正科技"
$ws.Range("A20").Value = "东方财富"
$ws.Range("A21").Value = "歌尔股份"

# Column C updates (rows 2-21)
$ws.Range("C2").Value = "岩山科技"
$ws.Range("C3").Value = "寒武纪"
$ws.Range("C4").Value = "华胜天成"
$ws.Range("C5").Value = "北方稀土"
$ws.Range("C6").Value = "天融信"
$ws.Range("C7").Value = "吉视传媒"
$ws.Range("C8").Value = "领益智造"
$ws.Range("C9").Value = "拓维信息"
$ws.Range("C10").Value = "万通发展"
$ws.Range("C11").Value = "剑桥科技"
$ws.Range("C12").Value = "利欧股份"
$ws.Range("C13").Value = "启明信息"
$ws.Range("C14").Value = "华银电力"
$ws.Range("C15").Value = "步步高"
$ws.Range("C16").Value = "光迅科技"
$ws.Range("C17").Value = "方正科技"
$ws.Range("C18").Value = "合力泰"
$ws.Range("C19").Value = "新易盛"
$ws.Range("C20").Value = "五粮液"
$ws.Range("C21").Value = "启明信息"
